$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header D1: "razon" -> "razoncierre"
$ws.Range("D1").Value = "razoncierre"

# Row 2: clear F2 and H2 entirely, set G2 to the new transaction id
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "AASIM23177L8VJBJ9"
$ws.Range("H2").ClearContents()

# Update the selection to D5
$ws.Range("D5").Select()
